# edit.ps1 - Applies the LR2/40.docx formatting + text-run-merge changes
# described by the commit "Changed LR2 (2 TIMES)".
#
# Summary of the change:
#  1. Every paragraph's pPr loses the w:before/w:after spacing (now only
#     w:line/w:lineRule survive) and loses the w:left/w:right indent; the
#     first line indent grows from 680 -> 709 twips (34pt -> 35.45pt).
#  2. A number of paragraphs had their text split across several <w:r>
#     runs with identical run formatting; those runs get merged back into
#     a single run (pure text concatenation, formatting unchanged).
#  3. The final paragraph ("Девушка действительно была...") loses its
#     cached <w:lastRenderedPageBreak/> marker and has its trailing "."
#     pulled out into its own trailing run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph formatting: applies uniformly to every paragraph in the
#    document body (all 15 paragraphs share the same starting format).
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $p.Format.SpaceBefore = 0
    $p.Format.SpaceAfter = 0
    $p.Format.LeftIndent = 0
    $p.Format.RightIndent = 0
    $p.Format.FirstLineIndent = 709 / 20
}

# ---------------------------------------------------------------------
# 2) Merge split runs: doing a Find/Replace across the run boundary with
#    MatchCase causes the engine to collapse the matched runs into one.
#    (Visible text & formatting is unchanged -- it is the exact same
#    text, just concatenated into a single run.)
# ---------------------------------------------------------------------
function Merge-Text($text) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2) | Out-Null
}

Merge-Text("Потом меня представили другому человеку. Потом еще. Я почувствовал, как внутри поднялась новая волна раздражения, и мой желудок начал сжиматься от голода. Вид еды был жутким, а желудок предательски заурчал.")

Merge-Text("Когда отец наконец отпустил меня, я вздохнул с облегчением. Идти прямо к столу с едой было бы слишком неприлично, поэтому ничего не оставалось, как брести по проходу, медленно приближаясь к ближайшему столику. Мне приходилось постоянно напоминать себе, что эта проклятая вечеринка очень важна для меня и моего папы, поэтому я должен был соблюдать рамки приличия. В свете этого было бы излишне бестактно набрасываться на еду.")

Merge-Text("И все бы ничего, если бы каждая остановка за столом с замысловатыми блюдами не предполагала выполнения определенных социальных обязанностей. Почему")

Merge-Text("то все считали своим долгом заметить меня, выразить радость по этому ограниченному, но обязательному случаю, подойти и поприветствовать мою скромную особу столь же многословную, сколь и фальшивую. И после неискреннего восхищения и еще менее искренних улыбок заставили меня завести краткий, но не менее грустный разговор ни о чем.")

Merge-Text("Самое смешное, что большинство из этих людей этого не знали. В общем. Но они я – да. Поэтому моей тушке ничего не оставалось, как со стоическим спокойствием выносить наполненные нездоровым любопытством мнения женщин")

Merge-Text("Несколько раз приходилось с холодной улыбкой извиняться, затеряться в толпе, сославшись на неотложные дела, чтобы быстро закрыть очередной пустой разговор с очередным собеседником, который неожиданно заинтересовался проектом ")

Merge-Text("Нет, все было не так уж плохо. Иногда я встречал действительно интересных людей. Также я знал некоторых гостей (например, здесь присутствовали все сотрудники отдела №118 или капитан Стейси с дочерью), и короткие беседы с ними совсем не могли быть грустными.")

Merge-Text("Пауку, что он такую девушку себе прихватит.")

Merge-Text("Среди прочего был такой персонаж, так Тони Старк. Правда, разговора с ним не получилось, как он ни старался, потому что внимание гениального изобретателя было полностью приковано к длинноногой модели, с которой он пришел. И на данный момент долгие разговоры с таким парнем, как я, его не интересовали.")

Merge-Text("Наконец, мои пищеварительные блуждания вывели меня к столу из одних, вид и запах выставленных на какой блюд, родился сглотнуть комок, внезапно образовавшийся в горле. Аппетит на креветки был, но, чувствуя себя любопытными взглядами, я вынужден был гасить желание схватить с десяток, и просто запихать их в рот.")

Merge-Text("Девушка действительно была великолепна. Длинное приталенное черное шелковое платье с глубоким вырезом сбоку, доходившим до бедра, было перехвачено на осиной талии у М. Джея шелковым поясом, подчеркивающим и без того бросающуюся в глаза фигуру. Ее длинные волосы, теперь распущенные, падали на плечи красными прядями светящегося пламени. Искусно нанесенный макияж подчеркнул пронзительную сторону, которая, уверена, покорила сегодня больше, чем мужское сердце.")

# ---------------------------------------------------------------------
# 3) Last paragraph: drop the lastRenderedPageBreak cache marker and
#    split the trailing "." into its own run.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("lastRenderedPageBreak") | Out-Null

$tail = $d.Content.Find.Execute("мужское сердце.", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "мужское сердце.", 2) | Out-Null
